$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.252.71"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.805.45"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'314.69"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.5274"
$ws.Range("E7").Value = "  +3.06%  "
$ws.Range("E8").Value = "  -3.11%  "
$ws.Range("D9").Value = "'0.08022"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "'41.42"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").Value = "'1.101"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "'6.331"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "'1.004"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").Value = "'20.62"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.808.12"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.334"
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").Value = "'92.21"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "'0.00001097"
$ws.Range("E18").Value = "  -3.48%  "
$ws.Range("D19").Value = "'0.06609"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("D23").Value = "28.315.93"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'11.16"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "'2.257"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").Value = "'160.93"
$ws.Range("E26").Value = "  +3.84%  "
$ws.Range("D27").Value = "'20.48"
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").Value = "2.011.32"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "'2.363"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").Value = "'123.43"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").Value = "'1.056"
$ws.Range("E32").Value = "  -3.98%  "
$ws.Range("D33").Value = "'3.685"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").Value = "'5.564"
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("D35").Value = "'0.07286"
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("D36").Value = "'12.51"
$ws.Range("E36").Value = "  +10.85%  "
$ws.Range("D37").Value = "'0.02319"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "'0.2160"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "'5.121"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").Value = "'8.669"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").Value = "'0.6211"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "'1.167"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").Value = "'1.369"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").Value = "'0.6034"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").Value = "'13.21"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").Value = "'3.772"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").Value = "'127.24"
$ws.Range("D48").Value = "'1.219"
$ws.Range("D49").Value = "'1.932"
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("D50").Value = "'0.06820"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").Value = "'73.29"
$ws.Range("E51").Value = "  -1.32%  "
